# Append two blank paragraphs followed by a new paragraph of text
# ("Batis lainnyha") at the end of the document, matching the paragraph
# mark formatting (lang=en-US) already used throughout the document.
#
# We build the new paragraphs as a raw WordprocessingML fragment and
# insert it with Range.InsertXML at a range collapsed to the very end of
# the document's content. Doing it as a single InsertXML call (rather
# than three separate Selection.TypeParagraph()/TypeText() calls) avoids
# synthesizing an empty <w:r> inside the blank paragraphs, matching how
# the paragraph mark alone carries the <w:rPr> (<w:lang w:val="en-US"/>)
# with no run, exactly like the rest of the document.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$blankPara = "<w:p $wNs><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr></w:p>"
$textPara  = "<w:p $wNs><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Batis lainnyha</w:t></w:r></w:p>"

$fragment = $blankPara + $blankPara + $textPara

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($fragment)
